$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (sandesh.achari): Book5 column now has a book checked out.
$ws.Range("H2").Value = "Inferno"

# Row 3 (sanket.achari): Book1/Book2 values updated; Book3 touched but left blank.
$ws.Range("D3").Value = "The Lost Symbol"
$ws.Range("E3").Value = "Open"
$ws.Range("F3").Style = "Normal"
